# Add a new "Concurrent Therapies" column (G) to the trials table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("G1").Value = "Adjuvant Therapy"

# Data rows (row number -> therapy text)
$ws.Range("G2").Value  = "Immunotherapy"
$ws.Range("G3").Value  = "Immunotherapy"
$ws.Range("G4").Value  = "Immunotherapy"
$ws.Range("G5").Value  = "Immunotherapy"
$ws.Range("G6").Value  = "Chemoradiation"
$ws.Range("G7").Value  = "Androgen deprivation therapy"
$ws.Range("G8").Value  = "Chemotherapy"
$ws.Range("G9").Value  = "Chemoradiation"
$ws.Range("G10").Value = "Chemoradiation"
$ws.Range("G11").Value = "Chemoradiation"
$ws.Range("G12").Value = "None"
$ws.Range("G13").Value = "None"
$ws.Range("G14").Value = "Chemoradiation"
$ws.Range("G15").Value = "Chemoradiation"
$ws.Range("G16").Value = "Chemo and/or Radiation therapy"
$ws.Range("G17").Value = "Chemoradiation"
$ws.Range("G18").Value = "Chemoradiation"
$ws.Range("G19").Value = "(SOC endocrine therapy and SGLT2i Therapy) or PI3K inhibition"
$ws.Range("G20").Value = "Stereotactic radiosurgery"
$ws.Range("G21").Value = "Chemoradiation"
$ws.Range("G22").Value = "N/A"
$ws.Range("G23").Value = "None"
$ws.Range("G24").Value = "None"

# Match formatting of the rest of the row (reuse the existing "data row" style)
# by copying formats from a representative already-styled cell in the same row.
# Column C is used as the format source because every row in C uses the common
# style (unlike column A, which has one special-cased cell at row 4).
$allRows = 1..24
foreach ($r in $allRows) {
    $ws.Range("C$r").Copy()
    $ws.Range("G$r").PasteSpecial(-4122)
}

# A few rows use the slightly different "no-family" font variant seen elsewhere
# in the sheet (B6, B7, B13) - replicate that formatting for the matching rows.
$ws.Range("B6").Copy()
$ws.Range("G5").PasteSpecial(-4122)

$ws.Range("B7").Copy()
$ws.Range("G7").PasteSpecial(-4122)

$ws.Range("B13").Copy()
$ws.Range("G13").PasteSpecial(-4122)

# Restore the active selection to match the post-edit cursor location.
$ws.Range("L19").Select()
